# Update "想去人数" (number of people interested) figures in both the
# "展览" and "全部类型" sheets to match the refreshed output snapshot.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 86
$ws1.Range("F4").Value = 267
$ws1.Range("F5").Value = 4
$ws1.Range("F6").Value = 10110
$ws1.Range("F7").Value = 330
$ws1.Range("F8").Value = 916
$ws1.Range("F9").Value = 1258
$ws1.Range("F10").Value = 6075
$ws1.Range("F12").Value = 417
$ws1.Range("F15").Value = 3107
$ws1.Range("F18").Value = 602
$ws1.Range("F20").Value = 19
$ws1.Range("F21").Value = 271
$ws1.Range("F23").Value = 1545

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 86
$ws4.Range("F5").Value = 267
$ws4.Range("F6").Value = 4
$ws4.Range("F7").Value = 10110
$ws4.Range("F8").Value = 330
$ws4.Range("F9").Value = 916
$ws4.Range("F10").Value = 1258
$ws4.Range("F11").Value = 6075
$ws4.Range("F13").Value = 417
$ws4.Range("F16").Value = 3107
$ws4.Range("F19").Value = 602
$ws4.Range("F21").Value = 19
$ws4.Range("F22").Value = 271
$ws4.Range("F24").Value = 1546
